$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date (column G)
$wsOverview.Range("G2").Value = "2016-10-27 10:28:32"
$wsOverview.Range("G3").Value = "2016-10-27 10:28:32"

# zh-cn sheet: Priority (column E) ht -> mt
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"

# zh-cn sheet: Correspond Handoff Datetime (column H)
$wsZhCn.Range("H2").Value = "2016-10-27 10:28:19"
$wsZhCn.Range("H3").Value = "2016-10-27 10:28:19"

# zh-cn sheet: Correspond Handback DateTime (column K)
$wsZhCn.Range("K2").Value = "2016-10-27 10:29:22"
$wsZhCn.Range("K3").Value = "2016-10-27 10:29:22"

# de-de sheet: Correspond Handback DateTime (column K)
$wsDeDe.Range("K2").Value = "2016-10-27 10:29:39"
$wsDeDe.Range("K3").Value = "2016-10-27 10:29:39"
